$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 280
$ws.Range("F7").Value = 13232
$ws.Range("F8").Value = 75
$ws.Range("F10").Value = 305
$ws.Range("F11").Value = 4791
$ws.Range("F13").Value = 3580
$ws.Range("F14").Value = 46
$ws.Range("F15").Value = 14
$ws.Range("F16").Value = 18
$ws.Range("F17").Value = 180
$ws.Range("F18").Value = 127
$ws.Range("F20").Value = 46
$ws.Range("F21").Value = 79
$ws.Range("F23").Value = 82
$ws.Range("F24").Value = 109
$ws.Range("F25").Value = 4373
$ws.Range("F27").Value = 1949
$ws.Range("F29").Value = 269
$ws.Range("F30").Value = 7067
$ws.Range("F31").Value = 22
$ws.Range("F32").Value = 169
$ws.Range("F33").Value = 2111
$ws.Range("F34").Value = 2062
$ws.Range("F35").Value = 1309
$ws.Range("F36").Value = 121
$ws.Range("F37").Value = 1098
$ws.Range("F38").Value = 16
$ws.Range("F40").Value = 234
$ws.Range("F41").Value = 228
$ws.Range("F43").Value = 8
$ws.Range("F44").Value = 152
$ws.Range("F45").Value = 1244
$ws.Range("F46").Value = 1869
$ws.Range("F47").Value = 80

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 132

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 484
$ws.Range("F3").Value = 660
$ws.Range("F4").Value = 42

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 484
$ws.Range("F6").Value = 660
$ws.Range("F7").Value = 280
$ws.Range("F8").Value = 13232
$ws.Range("F10").Value = 305
$ws.Range("F11").Value = 4792
$ws.Range("F12").Value = 3580
$ws.Range("F13").Value = 46
$ws.Range("F14").Value = 18
$ws.Range("F15").Value = 180
$ws.Range("F16").Value = 127
$ws.Range("F17").Value = 46
$ws.Range("F18").Value = 79
$ws.Range("F21").Value = 82
$ws.Range("F23").Value = 109
$ws.Range("F24").Value = 4373
$ws.Range("F26").Value = 1949
$ws.Range("F28").Value = 269
$ws.Range("F29").Value = 7067
$ws.Range("F31").Value = 22
$ws.Range("F32").Value = 169
$ws.Range("F33").Value = 2111
$ws.Range("F34").Value = 2062
$ws.Range("F35").Value = 1309
$ws.Range("F36").Value = 121
$ws.Range("F37").Value = 1098
$ws.Range("F39").Value = 234
$ws.Range("F40").Value = 228
$ws.Range("F42").Value = 152
$ws.Range("F44").Value = 1244
$ws.Range("F45").Value = 1869
$ws.Range("F46").Value = 80
